# Update "想去人数" (interest count) values in column F across all sheets
# per the commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 2397
$ws.Range("F9").Value = 137
$ws.Range("F11").Value = 1048
$ws.Range("F12").Value = 52
$ws.Range("F13").Value = 249
$ws.Range("F14").Value = 319
$ws.Range("F15").Value = 2123
$ws.Range("F16").Value = 1136
$ws.Range("F17").Value = 1032
$ws.Range("F18").Value = 818
$ws.Range("F19").Value = 84
$ws.Range("F20").Value = 825
$ws.Range("F21").Value = 1426
$ws.Range("F22").Value = 652
$ws.Range("F23").Value = 1664
$ws.Range("F24").Value = 33
$ws.Range("F25").Value = 314
$ws.Range("F26").Value = 56
$ws.Range("F27").Value = 90
$ws.Range("F29").Value = 2586

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 57
$ws.Range("F18").Value = 7
$ws.Range("F25").Value = 85
$ws.Range("F27").Value = 39
$ws.Range("F29").Value = 171
$ws.Range("F39").Value = 315
$ws.Range("F47").Value = 294

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1659
$ws.Range("F7").Value = 708
$ws.Range("F8").Value = 2478
$ws.Range("F9").Value = 9522
$ws.Range("F15").Value = 345
$ws.Range("F16").Value = 2735
$ws.Range("F17").Value = 341
$ws.Range("F18").Value = 643

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 708
$ws.Range("F11").Value = 2735
$ws.Range("F12").Value = 341
$ws.Range("F13").Value = 137
$ws.Range("F14").Value = 1048
$ws.Range("F16").Value = 643
$ws.Range("F22").Value = 249
$ws.Range("F23").Value = 1136
$ws.Range("F24").Value = 1032
$ws.Range("F25").Value = 818
$ws.Range("F26").Value = 84
$ws.Range("F27").Value = 825
$ws.Range("F32").Value = 652
$ws.Range("F35").Value = 1664
$ws.Range("F36").Value = 314
$ws.Range("F37").Value = 39
$ws.Range("F39").Value = 171
$ws.Range("F44").Value = 315
$ws.Range("F45").Value = 2586
$ws.Range("F48").Value = 294
